$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("DataSet")

# Leave a selection on DataSet matching where the user had scrolled to
$ds.Activate()
$ds.Range("A27").Select() | Out-Null

# Add the new sheet right after DataSet
$ws = $wb.Worksheets.Add($null, $ds)
$ws.Name = "Track_My_Order"

# --- Data row (row 2) kicked off first so new shared strings intern in
#     the same order as the authored workbook (GuestUserOrderdetails
#     before OrderID/prod order/Billinglastname/BillingEmail) ---
$ws.Range("A2").Value = "GuestUserOrderdetails"

# --- Header row (row 1) : same layout/style as DataSet's header row ---
$ws.Range("A1").Value = "DataSet"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Prod UserName"
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Confirm Password"
$ws.Range("F1").Value = "FirstName"
$ws.Range("G1").Value = "LastName"
$ws.Range("H1").Value = "OrderID"
$ws.Range("I1").Value = "prod order"
$ws.Range("J1").Value = "Billinglastname"
$ws.Range("K1").Value = "BillingEmail"
$ws.Range("L1").Value = "Prod Email"
$ws.Range("M1").Value = "methods"
$ws.Range("N1").Value = "Street"
$ws.Range("O1").Value = "City"
$ws.Range("P1").Value = "Country"
$ws.Range("Q1").Value = "Region"
$ws.Range("R1").Value = "postcode"
$ws.Range("S1").Value = "phone"
$ws.Range("T1").Value = "OTP Number"
$ws.Range("U1").Value = "cardNumber"
$ws.Range("V1").Value = "ExpMonthYear"
$ws.Range("W1").Value = "cvv"
$ws.Range("X1").Value = "Products"
$ws.Range("Y1").Value = "Quantity"
$ws.Range("Z1").Value = "Colorproduct"
$ws.Range("AA1").Value = "Color"
$ws.Range("AB1").Value = "PLP Color"
$ws.Range("AC1").Value = "Discountcode"
$ws.Range("AD1").Value = "prodDiscountcode"
$ws.Range("AE1").Value = "DOB"
$ws.Range("AF1").Value = "Links"
$ws.Range("AG1").Value = "Sort"
$ws.Range("AH1").Value = "message"

$ws.Range("A1:AH1").Interior.Color = 65535

# --- Data row (row 2) continued ---
$ws.Range("H2").Value = "DRYUSSTG3000001743"
$ws.Range("J2").Value = "Qa"
$ws.Range("K2").Value = "avayugundla@helenoftroy.com"
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:avayugundla@helenoftroy.com") | Out-Null

# --- Column widths matching the header content ---
$ws.Columns.Item(1).ColumnWidth = 18.830729166666668
$ws.Columns.Item(8).ColumnWidth = 19.053385416666668
$ws.Columns.Item(9).ColumnWidth = 8.830729166666666
$ws.Columns.Item(10).ColumnWidth = 12.498697916666666
$ws.Columns.Item(11).ColumnWidth = 26.166666666666668

# --- View state: new sheet active, selection on H4 ---
$ws.Range("H4").Select() | Out-Null
